$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.616400480270386
$ws.Range("B1").Value = 2.695095777511597
$ws.Range("C1").Value = 3.05919623374939
$ws.Range("D1").Value = 3.424724578857422
$ws.Range("E1").Value = 2.068240642547607
